# Minor adjustments to presentation: reposition four pictures on slide 5.
# PowerPoint's Shape.Left/.Top are Single-precision points; a small
# (0.5 EMU) epsilon is added before the EMU->point division so that the
# float32 round-trip back to EMU lands exactly on the target value
# instead of being floored one EMU short.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Picture 4 (rId2): off 4411996,577850 -> 4545955,825605 (EMU)
$shPic4 = $s.Shapes.Item(10)
$shPic4.Left = (4545955 + 0.5) / 12700
$shPic4.Top  = (825605 + 0.5) / 12700

# Picture 7 (rId3): off 8574192,4483395 -> 4680313,3909732 (EMU)
$shPic7 = $s.Shapes.Item(11)
$shPic7.Left = (4680313 + 0.5) / 12700
$shPic7.Top  = (3909732 + 0.5) / 12700

# Picture 12 (rId4): off 5171735,3737549 -> 7440518,4213451 (EMU)
$shPic12 = $s.Shapes.Item(12)
$shPic12.Left = (7440518 + 0.5) / 12700
$shPic12.Top  = (4213451 + 0.5) / 12700

# Picture 8 (rId5): off 9194217,577850 -> 9291468,587644 (EMU)
$shPic8 = $s.Shapes.Item(13)
$shPic8.Left = (9291468 + 0.5) / 12700
$shPic8.Top  = (587644 + 0.5) / 12700
